# Commit: "removes the term 'slave' from the project where possible."
#
# The Modbus help sheet contains two shared-string cells that used the word
# "slave". Replace them with the updated wording, preserving everything
# else (cell formatting / row heights are carried over automatically since
# only the .Value changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modbus")

# A4: "... Inputs with the slave id set to 0 are turned off."
#  -> "... Inputs with the device id set to 0 are turned off."
$ws.Range("A4").Value = "The inputs 1+2 configure the MODBUS device, inputs 3+4 configure the MODBUS_34 device and so on.`nInputs with the device id set to 0 are turned off."

# A9: "... connected PID slave using ... turn the PID slave on and off ..."
#  -> "... connected PID using ... turn the PID  on and off ..."
$ws.Range("A9").Value = "The PID Control dialog can operate a connected PID using the given PID registers to set the p-i-d parameters and the set value (SV). MODBUS commands can be specified to turn the PID  on and off from that PID Control dialog. See the help page in the Events Dialog for documentation of available MODBUS write commands."

# The saved workbook (re-exported by whatever app made this edit upstream)
# also shows the row heights of A4/A5 settling to auto-fitted values and the
# active selection moving to A13 -- replicate those observable, low-risk
# cosmetic deltas too.
$ws.Rows("4:4").RowHeight = 22.7
$ws.Rows("5:5").RowHeight = 69.3

$ws.Range("A13").Select()
